$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "https://www.capitalonecareers.com"

for ($row = 2; $row -le 25; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $relative = $cell.Value2
    $full = $prefix + $relative
    $ws.Hyperlinks.Add($cell, $full)
    $cell.Value2 = $full
}
